$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 2 values ---
$ws.Range("D2").Value = 1000

$ws.Range("K2").Value = 0.9902999997138977
$ws.Range("L2").Value = 0.9814000129699707
$ws.Range("M2").Value = 432.465
$ws.Range("N2").Value = 0.0023
$ws.Range("O2").Value = 0.0023
$ws.Range("P2").Value = 54
$ws.Range("Q2").Value = 8.008599999999999
$ws.Range("R2").Value = 0.9901999831199646
$ws.Range("S2").Value = 0.991100013256073
$ws.Range("T2").Value = 0.991100013256073

# --- Add new row 3 ---
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 0.0001
$ws.Range("C3").Value = 3
$ws.Range("D3").Value = 1000
$ws.Range("E3").Value = 200
$ws.Range("F3").Value = 4
$ws.Range("G3").Value = "('tanh', 'relu')"
$ws.Range("H3").Value = 100
$ws.Range("I3").Value = 0.0001
$ws.Range("J3").Value = 2
$ws.Range("K3").Value = 0.9965000152587891
$ws.Range("L3").Value = 0.9916999936103821
$ws.Range("M3").Value = 541.039
$ws.Range("N3").Value = 0.0018
$ws.Range("O3").Value = 0.0018
$ws.Range("P3").Value = 66
$ws.Range("Q3").Value = 8.1976
$ws.Range("R3").Value = 0.9957000017166138
$ws.Range("S3").Value = 0.995199978351593
$ws.Range("T3").Value = 0.995199978351593

# --- Add new row 4 ---
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 0.0001
$ws.Range("C4").Value = 3
$ws.Range("D4").Value = 1000
$ws.Range("E4").Value = 200
$ws.Range("F4").Value = 4
$ws.Range("G4").Value = "('tanh', 'relu')"
$ws.Range("H4").Value = 100
$ws.Range("I4").Value = 0.0001
$ws.Range("J4").Value = 3
$ws.Range("K4").Value = 0.9980000257492065
$ws.Range("L4").Value = 0.9959999918937683
$ws.Range("M4").Value = 631.567
$ws.Range("N4").Value = 0.0016
$ws.Range("O4").Value = 0.0016
$ws.Range("P4").Value = 80
$ws.Range("Q4").Value = 7.8946
$ws.Range("R4").Value = 0.9977999925613403
$ws.Range("S4").Value = 0.9980000257492065
$ws.Range("T4").Value = 0.9980000257492065

# Copy style (bold/border/centered) from A2 to A3 and A4 so the new rows match formatting
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A3:A4").PasteSpecial(-4122) | Out-Null
